$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the manufacturer value in N2 (was "GM Manufacturer")
$ws.Range("N2").Value = ""

# Update category text in AL2 and AM2: replace the separator before the
# last segment ("> Subcategory") with "| Subcategory"
$ws.Range("AL2").Value = "Tovary a kategórie > GM Category | Subcategory"
$ws.Range("AM2").Value = "Tovary a kategórie > GM Category | Subcategory"

Write-Output "edit applied"
